{"js": "const pairs = [\n  [\"602\u00d77=4214\", \"407\u00d75=2035\"],\n  [\"704\u00d76=4224\", \"811\u00d73=2433\"],\n  [\"932\u00d72=1864\", \"550\u00d78=4400\"],\n  [\"843\u00d73=2529\", \"213\u00d76=1278\"],\n  [\"236\u00d75=1180\", \"918\u00d79=8262\"],\n  [\"164\u00d79=1476\", \"274\u00d74=1096\"],\n  [\"175\u00d72=350\", \"673\u00d72=1346\"],\n  [\"138\u00d73=414\", \"708\u00d78=5664\"],\n  [\"995\u00d74=3980\", \"364\u00d78=2912\"],\n  [\"142\u00d77=994\", \"312\u00d77=2184\"],\n  [\"733\u00d76=4398\", \"836\u00d76=5016\"],\n  [\"873\u00d74=3492\", \"477\u00d79=4293\"],\n  [\"787\u00d74=3148\", \"755\u00d76=4530\"],\n  [\"495\u00d74=1980\", \"114\u00d72=228\"],\n  [\"950\u00d77=6650\", \"840\u00d77=5880\"],\n  [\"238\u00d73=714\", \"506\u00d73=1518\"],\n  [\"965\u00d78=7720\", \"369\u00d74=1476\"],\n  [\"439\u00d76=2634\", \"501\u00d77=3507\"],\n  [\"987\u00d74=3948\", \"138\u00d75=690\"],\n  [\"232\u00d78=1856\", \"319\u00d74=1276\"],\n  [\"155\u00d73=465\", \"760\u00d79=6840\"],\n  [\"393\u00d77=2751\", \"926\u00d73=2778\"],\n  [\"282\u00d76=1692\", \"224\u00d75=1120\"],\n  [\"213\u00d75=1065\", \"759\u00d77=5313\"],\n  [\"196\u00d72=392\", \"534\u00d75=2670\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"602\u00d77=4214\", \"407\u00d75=2035\"),\n  @(\"704\u00d76=4224\", \"811\u00d73=2433\"),\n  @(\"932\u00d72=1864\", \"550\u00d78=4400\"),\n  @(\"843\u00d73=2529\", \"213\u00d76=1278\"),\n  @(\"236\u00d75=1180\", \"918\u00d79=8262\"),\n  @(\"164\u00d79=1476\", \"274\u00d74=1096\"),\n  @(\"175\u00d72=350\", \"673\u00d72=1346\"),\n  @(\"138\u00d73=414\", \"708\u00d78=5664\"),\n  @(\"995\u00d74=3980\", \"364\u00d78=2912\"),\n  @(\"142\u00d77=994\", \"312\u00d77=2184\"),\n  @(\"733\u00d76=4398\", \"836\u00d76=5016\"),\n  @(\"873\u00d74=3492\", \"477\u00d79=4293\"),\n  @(\"787\u00d74=3148\", \"755\u00d76=4530\"),\n  @(\"495\u00d74=1980\", \"114\u00d72=228\"),\n  @(\"950\u00d77=6650\", \"840\u00d77=5880\"),\n  @(\"238\u00d73=714\", \"506\u00d73=1518\"),\n  @(\"965\u00d78=7720\", \"369\u00d74=1476\"),\n  @(\"439\u00d76=2634\", \"501\u00d77=3507\"),\n  @(\"987\u00d74=3948\", \"138\u00d75=690\"),\n  @(\"232\u00d78=1856\", \"319\u00d74=1276\"),\n  @(\"155\u00d73=465\", \"760\u00d79=6840\"),\n  @(\"393\u00d77=2751\", \"926\u00d73=2778\"),\n  @(\"282\u00d76=1692\", \"224\u00d75=1120\"),\n  @(\"213\u00d75=1065\", \"759\u00d77=5313\"),\n  @(\"196\u00d72=392\", \"534\u00d75=2670\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $found = $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n  if (-not $found) {\n    throw \"Replace failed: could not find '\" + $pair[0] + \"'\"\n  }\n}\n"}
